$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 2
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 4
$ws.Range("D1").Formula = "=SOMA(A1:C1)"
